$wb = $excel.ActiveWorkbook

# --- 1. Fix cell text: "Loop" -> "Looping" (appears once per sheet, same
#        logical cell: class=Elementary / activityReference=TestItem_Loop:0 / name=Loop) ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value() -eq "Loop") {
                $cell.Value = "Looping"
            }
        }
    }
}

# --- 2. Rename the sheets to their final names (do this before reordering,
#        since Move() re-seats old variable handles by position) ---
$wb.Worksheets.Item("TestItem_SeqWithLoop").Name = "SequenceWithLoop"
$wb.Worksheets.Item("TestItem_StartWithLoop").Name = "StartWithLoop"

# --- 3. Reorder tabs: StartWithLoop first, SequenceWithLoop second ---
$wb.Worksheets.Item("StartWithLoop").Move($wb.Worksheets.Item("SequenceWithLoop"))

# --- 4. Restore selections / active sheet as left by the author ---
$startSheet = $wb.Worksheets.Item("StartWithLoop")
$seqSheet = $wb.Worksheets.Item("SequenceWithLoop")

[void]$startSheet.Range("C5").Select()
[void]$seqSheet.Range("C6").Select()
[void]$seqSheet.Activate()
